$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FRED codes")
$ws2 = $wb.Worksheets.Item("Portfolio")

# Remove the three "Canada" entries from the "FRED codes" sheet (rows 11-13:
# 3-Month rate, immediate rate, 10-Year bond yield for Canada) and shift the
# remaining rows up.
$ws1.Rows("11:13").Delete()

# Remove the "Canada" entry from the "Portfolio" sheet (row 5) and shift the
# remaining rows up.
$ws2.Rows("5:5").Delete()

# The hidden AutoFilter-database named ranges need to shrink to match the new
# (smaller) data extents on each sheet.
$wb.Names.Item("FRED codes!_FilterDatabase").RefersTo = "='FRED codes'!`$A`$1:`$F`$61"
$wb.Names.Item("Portfolio!_FilterDatabase").RefersTo = "=Portfolio!`$A`$1:`$D`$21"

# Restore a sensible selection on "FRED codes" ...
$ws1.Activate()
$ws1.Range("E1").Select()

# ... then make "Portfolio" the active sheet/tab, with the entire row 5
# selected (the row that moved into the old Canada slot).
$ws2.Activate()
$ws2.Range("A5:XFD5").Select()
